# Update need_to_buy.xlsx values to reflect the refreshed R calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ column letter = new value }
$changes = @{
    2  = @{ B = 11576.322163638;  C = 11501.1394459734; E = 7487.70056213383; F = 8.04083367113405 }
    3  = @{ B = 11965.9310003792; C = 11381.3498318679; E = 7622.78504964094; F = 288.678120062866 }
    4  = @{ B = 11786.0520238695; C = 11433.6144150008; E = 7915.94902436248; F = 303.070976640136 }
    5  = @{ B = 11958.1035091571; C = 11003.417607155;  E = 8053.17566448143; F = 290.863886318185 }
    6  = @{ B = 4853.60248239427; C = 7881.48987791606; E = 7950.08737376963; F = 156.488218820237 }
    7  = @{ B = 5097.80826294459; C = 7920.00923285382; E = 7959.97255022592; F = 253.671740961656 }
    9  = @{ C = 10821.9595269282; F = 368.566674314826 }
    10 = @{ C = 10771.6945650184; F = 366.472300901915 }
    11 = @{ C = 10883.6423363379; F = 371.136791373563 }
    12 = @{ C = 10333.2663600624; F = 348.204459028749 }
    13 = @{ C = 7328.55065814478; F = 207.451638146692 }
    14 = @{ C = 7436.54841837552; F = 211.59015638758 }
    15 = @{ C = 10856.5243006978; F = 350.159195785845 }
}

foreach ($row in $changes.Keys) {
    foreach ($col in $changes[$row].Keys) {
        $address = "$col$row"
        $ws.Range($address).Value = $changes[$row][$col]
    }
}

$wb.Save()
